{"js": "// Load all paragraphs once so we can locate the four \"Methods used:\" lines\n// and the Conclusion heading (only needed for context / sanity, the actual\n// bookmark id bookkeeping is handled automatically by the host when a new\n// bookmark is inserted).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Small helper: find `searchText` inside paragraph `para` and replace the\n// whole matched (possibly multi-run) span with `replacementText`.\nasync function replaceInParagraph(para, searchText, replacementText) {\n  const range = para.getRange(Word.RangeLocation.content);\n  const results = range.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 1) Methods used: Literature study, Brainstorm, Design pattern research\n//    -> Literature study, Expert interview, Design pattern research,\n//       Document analysis, Brainstorm, Problem analysis\n// ---------------------------------------------------------------------\nawait replaceInParagraph(\n  paragraphs.items[57],\n  \"Literature study, Brainstorm, Design pattern research\",\n  \"Literature study, Expert interview, Design pattern research, Document analysis, Brainstorm, Problem analysis\"\n);\n\n// ---------------------------------------------------------------------\n// 2) Methods used: Literature study, Problem analysis, Pitch\n//    -> Literature study, Expert interview, Document analysis, Problem\n//       analysis, Pitch\n// ---------------------------------------------------------------------\nawait replaceInParagraph(\n  paragraphs.items[61],\n  \"Literature study, Problem analysis, Pitch\",\n  \"Literature study, Expert interview, Document analysis, Problem analysis, Pitch\"\n);\n\n// ---------------------------------------------------------------------\n// 3) Methods used: Literature study, Available product analysis, Expert\n//    Interview, Best good and bad practices\n//    -> Literature study, Expert interview, Available product analysis,\n//       Prototyping, Problem analysis, System test\n// ---------------------------------------------------------------------\nawait replaceInParagraph(\n  paragraphs.items[65],\n  \"Literature study, Available product analysis, Expert Interview, Best good and bad practices\",\n  \"Literature study, Expert interview, Available product analysis, Prototyping, Problem analysis, System test\"\n);\n\n// ---------------------------------------------------------------------\n// 4) Methods used: Literature study, Problem analysis, Best good and bad\n//    practices, Expert Interview, Non-functional test, Security test\n//    -> \"Best good and bad practices\" becomes \"Component Test\" AND the\n//    whole methods list (everything after \"Methods used: \") gets wrapped\n//    in a new bookmark named \"_Hlk161407544\". Any existing bookmark whose\n//    numeric id collides with the new one (here the Conclusion heading's\n//    _Toc161166472) is renumbered automatically by the host.\n// ---------------------------------------------------------------------\nconst p4 = paragraphs.items[69];\n\nconst bmRange = p4.getRange(Word.RangeLocation.content);\nconst bmResults = bmRange.search(\n  \"Literature study, Problem analysis, Best good and bad practices, Expert Interview, Non-functional test, Security test\",\n  { matchCase: true }\n);\nbmResults.load(\"text\");\nawait context.sync();\nbmResults.items[0].insertBookmark(\"_Hlk161407544\");\nawait context.sync();\n\nawait replaceInParagraph(p4, \"Best good and bad practices, \", \"Component Test, \");\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# Helper: replace `oldText` with `newText` inside a given paragraph's range\n# ---------------------------------------------------------------------------\nfunction Replace-InParagraph($paraIndex, $oldText, $newText) {\n    $p = $d.Paragraphs.Item($paraIndex)\n    $r = $p.Range\n    $find = $r.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null,$null,$null,$null,$null,$null,$null,$null,$null,$null,2) | Out-Null\n}\n\n# ---------------------------------------------------------------------------\n# 1) \"Methods used:\" paragraph for the architecture-type sub-question\n#    Literature study, Brainstorm, Design pattern research\n#    -> Literature study, Expert interview, Design pattern research,\n#       Document analysis, Brainstorm, Problem analysis\n# ---------------------------------------------------------------------------\nReplace-InParagraph 58 `\n    \"Literature study, Brainstorm, Design pattern research\" `\n    \"Literature study, Expert interview, Design pattern research, Document analysis, Brainstorm, Problem analysis\"\n\n# ---------------------------------------------------------------------------\n# 2) \"Methods used:\" paragraph for the \"developing\" technologies sub-question\n#    Literature study, Problem analysis, Pitch\n#    -> Literature study, Expert interview, Document analysis, Problem analysis, Pitch\n# ---------------------------------------------------------------------------\nReplace-InParagraph 62 `\n    \"Literature study, Problem analysis, Pitch\" `\n    \"Literature study, Expert interview, Document analysis, Problem analysis, Pitch\"\n\n# ---------------------------------------------------------------------------\n# 3) \"Methods used:\" paragraph for the \"deploying\" technologies sub-question\n#    Literature study, Available product analysis, Expert Interview, Best good and bad practices\n#    -> Literature study, Expert interview, Available product analysis, Prototyping,\n#       Problem analysis, System test\n# ---------------------------------------------------------------------------\nReplace-InParagraph 66 `\n    \"Literature study, Available product analysis, Expert Interview, Best good and bad practices\" `\n    \"Literature study, Expert interview, Available product analysis, Prototyping, Problem analysis, System test\"\n\n# ---------------------------------------------------------------------------\n# 4) \"Methods used:\" paragraph for the \"testing\" sub-question\n#    Best good and bad practices -> Component Test\n#    and the whole methods list gets wrapped in a new bookmark \"_Hlk161407544\"\n#    (Word automatically renumbers any colliding bookmark ids, e.g. the\n#    Conclusion heading's _Toc161166472 bookmark).\n# ---------------------------------------------------------------------------\n$p4 = $d.Paragraphs.Item(70)\n\n# Build a range spanning from the start of the methods list (\"Literature\n# study...\") to the end of the paragraph text (excluding the paragraph mark).\n$bmRange = $p4.Range\n$bmRange.MoveEnd(1, -1) | Out-Null\n$bmFind = $bmRange.Find\n$bmFind.Text = \"Literature study, Problem analysis, Best good and bad practices, Expert Interview, Non-functional test, Security test\"\n$bmFind.Execute() | Out-Null\n$d.Bookmarks.Add(\"_Hlk161407544\", $bmRange) | Out-Null\n\nReplace-InParagraph 70 `\n    \"Best good and bad practices, \" `\n    \"Component Test, \"\n"}
